# Update attendance/ticket-count figures (column F) on the "展览" and
# "全部类型" worksheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 8313
$ws1.Range("F5").Value = 6061
$ws1.Range("F6").Value = 521
$ws1.Range("F10").Value = 314
$ws1.Range("F11").Value = 996

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 8313
$ws4.Range("F5").Value = 6061
$ws4.Range("F6").Value = 521
$ws4.Range("F10").Value = 314
$ws4.Range("F15").Value = 996
